$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 3, 11.18077166666667, 33.542315, 0.0566654123397129, 0.05684167636225704, 59.81642030112778, 538.3477827101501, 0.002350620346357435, 0.002357932210519635),
    @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 3, 12.200174, 36.600522, 0.06183185838481134, 0.06202419320830027, 65.27015822231333, 587.43142400082, 0.002564937205452364, 0.002572915725871411),
    @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 3, 84.38991299999999, 253.169739, 0.4276976008475411, 0.4290280014648685, 451.48068985551, 4063.32620869959, 0.01774194594426179, 0.01779713422633318),
    @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 3, 87.70567066666668, 263.117012, 0.4445022348210965, 0.4458849081081047, 469.2197833741912, 4222.978050367721, 0.01843904339578152, 0.01849640007645511),
    @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 2, 1.8355735, 3.671147, 0.009302893606838107, 0.006221220856469531, 9.820201972011667, 58.92121183207, 0.0003859068537458904, 0.000258071506419653),
    @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 3, 11.18077166666667, 33.542315, 0.0566654123397129, 0.05684167636225704, 15.21893096951667, 136.97037872565, 0.0005980620138527527, 0.0005999223517897742),
    @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 3, 12.200174, 36.600522, 0.06183185838481134, 0.06202419320830027, 16.60651084358, 149.45859759222, 0.0006525900760094222, 0.0006546200295052194),
    @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 3, 84.38991299999999, 253.169739, 0.4276976008475411, 0.4290280014648685, 114.86901787821, 1033.82116090389, 0.004514035598106921, 0.004528077004311816),
    @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 3, 87.70567066666668, 263.117012, 0.4445022348210965, 0.4458849081081047, 119.3823277413467, 1074.44094967212, 0.004691396228186365, 0.004705989334216742),
    @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 2, 1.8355735, 3.671147, 0.009302893606838107, 0.006221220856469531, 2.498527580995, 14.99116548597, 0.00009818524308635935, 0.00006566043941826835),
    @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 3, 11.18077166666667, 33.542315, 0.0566654123397129, 0.05684167636225704, 503.0121017502079, 4527.108915751871, 0.01976698830999357, 0.01982847570996543),
    @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 3, 12.200174, 36.600522, 0.06183185838481134, 0.06202419320830027, 548.8740266250173, 4939.866239625157, 0.02156923547207944, 0.02163632896086794),
    @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 3, 84.38991299999999, 253.169739, 0.4276976008475411, 0.4290280014648685, 3796.620552694159, 34169.58497424742, 0.1491967167816867, 0.1496608096447662),
    @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 3, 87.70567066666668, 263.117012, 0.4445022348210965, 0.4458849081081047, 3945.793282674576, 35512.13954407119, 0.15505879365704, 0.1555411211575791),
    @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 2, 1.8355735, 3.671147, 0.009302893606838107, 0.006221220856469531, 82.58067615356768, 495.4840569214061, 0.003245192818381854, 0.002170191566002897),
    @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 3, 11.18077166666667, 33.542315, 0.0566654123397129, 0.05684167636225704, 863.9217387720488, 7775.29564894844, 0.03394974166950915, 0.0340553460899822),
    @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 3, 12.200174, 36.600522, 0.06183185838481134, 0.06202419320830027, 942.6894537900745, 8484.205084110672, 0.03704509563127012, 0.0371603284920557),
    @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 3, 84.38991299999999, 253.169739, 0.4276976008475411, 0.4290280014648685, 6520.684130518295, 58686.15717466467, 0.2562449025234858, 0.2570419805894573),
    @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 3, 87.70567066666668, 263.117012, 0.4445022348210965, 0.4458849081081047, 6776.887835784324, 60991.99052205892, 0.2663130015400887, 0.2671413975398538),
    @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 2, 1.8355735, 3.671147, 0.009302893606838107, 0.006221220856469531, 141.8320574859453, 850.992344915672, 0.005573608691624004, 0.003727297344628714)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $vals = $data[$i]
    $ws.Cells.Item($r, 5).Value = $vals[0]
    $ws.Cells.Item($r, 7).Value = $vals[1]
    $ws.Cells.Item($r, 8).Value = $vals[2]
    $ws.Cells.Item($r, 9).Value = $vals[3]
    $ws.Cells.Item($r, 10).Value = $vals[4]
    $ws.Cells.Item($r, 11).Value = $vals[5]
    $ws.Cells.Item($r, 13).Value = $vals[6]
    $ws.Cells.Item($r, 14).Value = $vals[7]
    $ws.Cells.Item($r, 15).Value = $vals[8]
    $ws.Cells.Item($r, 16).Value = $vals[9]
    $ws.Cells.Item($r, 17).Value = $vals[10]
    $ws.Cells.Item($r, 18).Value = $vals[11]
    $ws.Cells.Item($r, 19).Value = $vals[12]
    $ws.Cells.Item($r, 20).Value = $vals[13]
}
